$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: 陳毅芸 / 清洗資料、存放資料 / 安裝mongoDB / 寫好初步code / 儲存VIX資料
# (new shared-string cells are populated F before E so the new unique strings land
#  in the shared string table in the same order as the target workbook: 儲存VIX資料(31), 寫好初步code(32))
$ws.Range("A8").Value = 12.2
$ws.Range("B8").Value = "陳毅芸"
$ws.Range("C8").Value = "清洗資料、存放資料"
$ws.Range("D8").Value = "安裝mongoDB"
$ws.Range("F8").Value = "儲存VIX資料"
$ws.Range("E8").Value = "寫好初步code"

# Row 9: 吳培瑜 / 爬蟲 用selenium爬CNBC中market和finance和amazon相關的新聞標題 / code完成 /
#        爬完cnbc從2007至今所有新聞標題 / 研究如何爬amazon相關更久遠以前的新聞
$ws.Range("A9").Value = 12.2
$ws.Range("B9").Value = "吳培瑜"
$ws.Range("C9").Value = "爬蟲 用selenium爬CNBC中market和finance和amazon相關的新聞標題"
$ws.Range("D9").Value = "code完成"
$ws.Range("E9").Value = "爬完cnbc從2007至今所有新聞標題"
$ws.Range("F9").Value = "研究如何爬amazon相關更久遠以前的新聞"

# Update the selection / active cell to match the saved view state (F9)
$ws.Range("F9").Select()
